# Auto-generated Excel COM-interop script applying numeric updates
# described by the Hyperion_Profits sheet diff (scheduled runner refresh).
$wb = $excel.ActiveWorkbook

# --- ALC row 17 ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 3901.647
$ws.Range("J17").Value = 4356.567
$ws.Range("L17").Value = 13069.701
$ws.Range("N17").Value = -13405.701

# --- ALC row 80 ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H80").Value = 972
$ws.Range("I80").Value = 839
$ws.Range("J80").Value = 1105
$ws.Range("K80").Value = 2517
$ws.Range("L80").Value = 3315
$ws.Range("M80").Value = -1519
$ws.Range("N80").Value = -5311

# --- ALC row 83 ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H83").Value = 972
$ws.Range("I83").Value = 839
$ws.Range("J83").Value = 1105
$ws.Range("K83").Value = 7551
$ws.Range("L83").Value = 9945
$ws.Range("M83").Value = -2559
$ws.Range("N83").Value = -19929

# --- ALC row 92 ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H92").Value = 864.0714
$ws.Range("I92").Value = 888.2
$ws.Range("K92").Value = 888.2
$ws.Range("M92").Value = 359.8

# --- ALC row 96 ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H96").Value = 786.8461
$ws.Range("I96").Value = 681.375
$ws.Range("J96").Value = 955.6
$ws.Range("K96").Value = 2044.125
$ws.Range("L96").Value = 2866.8
$ws.Range("M96").Value = -671.125
$ws.Range("N96").Value = -5612.8

# --- ALC row 112 ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H112").Value = 12480
$ws.Range("J112").Value = 13890
$ws.Range("L112").Value = 41670
$ws.Range("N112").Value = -43886

# --- ALC row 137 ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 56665.273
$ws.Range("I137").Value = 91201.75
$ws.Range("J137").Value = 3532.2307
$ws.Range("K137").Value = 273605.25
$ws.Range("L137").Value = 10596.6921
$ws.Range("M137").Value = -271055.25
$ws.Range("N137").Value = -15696.6921

# --- ALC row 138 ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 3228.1177
$ws.Range("J138").Value = 4054.125
$ws.Range("L138").Value = 12162.375
$ws.Range("N138").Value = -22442.375

# --- ARM row 2 ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2456
$ws.Range("I2").Value = 1276.3334
$ws.Range("J2").Value = 3340.75
$ws.Range("K2").Value = 1276.3334
$ws.Range("L2").Value = 3340.75
$ws.Range("M2").Value = -1163.3334
$ws.Range("N2").Value = -3566.75

# --- ARM row 61 ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 1811.6296
$ws.Range("J61").Value = 1852.2
$ws.Range("L61").Value = 1852.2
$ws.Range("N61").Value = -2276.2

# --- ARM row 102 ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H102").Value = 3820.8572
$ws.Range("I102").Value = 3076.1875
$ws.Range("K102").Value = 3076.1875
$ws.Range("M102").Value = -1454.1875

# --- ARM row 116 ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H116").Value = 2456
$ws.Range("I116").Value = 1276.3334
$ws.Range("J116").Value = 3340.75
$ws.Range("K116").Value = 1276.3334
$ws.Range("L116").Value = 3340.75
$ws.Range("M116").Value = 1017.6666
$ws.Range("N116").Value = -7928.75

# --- ARM row 122 ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 20206042
$ws.Range("I122").Value = 37039990
$ws.Range("J122").Value = 5298.4
$ws.Range("K122").Value = 111119970
$ws.Range("L122").Value = 15895.2
$ws.Range("M122").Value = -111117520
$ws.Range("N122").Value = -20795.2

# --- ARM row 136 ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 1811.6296
$ws.Range("J136").Value = 1852.2
$ws.Range("L136").Value = 5556.6
$ws.Range("N136").Value = -10656.6

# --- BSM row 3 ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2456
$ws.Range("I3").Value = 1276.3334
$ws.Range("J3").Value = 3340.75
$ws.Range("K3").Value = 1276.3334
$ws.Range("L3").Value = 3340.75
$ws.Range("M3").Value = -1162.3334
$ws.Range("N3").Value = -3568.75

# --- BSM row 94 ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 5554.1577
$ws.Range("I94").Value = 1053.3334
$ws.Range("J94").Value = 13269.857
$ws.Range("K94").Value = 1053.3334
$ws.Range("L94").Value = 13269.857
$ws.Range("M94").Value = -602.3334
$ws.Range("N94").Value = -14171.857

# --- BSM row 99 ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 3790.6667
$ws.Range("I99").Value = 2856.2856
$ws.Range("J99").Value = 5098.8
$ws.Range("K99").Value = 2856.2856
$ws.Range("L99").Value = 5098.8
$ws.Range("M99").Value = -1358.2856
$ws.Range("N99").Value = -8094.8

# --- CRP row 31 ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 14979.732
$ws.Range("I31").Value = 1819.4445
$ws.Range("J31").Value = 16517.947
$ws.Range("K31").Value = 1819.4445
$ws.Range("L31").Value = 16517.947
$ws.Range("M31").Value = -1524.4445
$ws.Range("N31").Value = -17107.947

# --- CRP row 34 ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 14979.732
$ws.Range("I34").Value = 1819.4445
$ws.Range("J34").Value = 16517.947
$ws.Range("K34").Value = 1819.4445
$ws.Range("L34").Value = 16517.947
$ws.Range("M34").Value = -1617.4445
$ws.Range("N34").Value = -16921.947

# --- CRP row 52 ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H52").Value = 0
$ws.Range("J52").Value = 0
$ws.Range("L52").Value = 0
$ws.Range("N52").ClearContents()

# --- CRP row 99 ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 3427.7693
$ws.Range("J99").Value = 4798
$ws.Range("L99").Value = 4798
$ws.Range("N99").Value = -7794

# --- CRP row 122 ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H122").Value = 2249.4333
$ws.Range("I122").Value = 1705.6086
$ws.Range("K122").Value = 5116.825800000001
$ws.Range("M122").Value = -2666.825800000001

# --- CRP row 126 ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H126").Value = 3427.7693
$ws.Range("J126").Value = 4798
$ws.Range("L126").Value = 14394
$ws.Range("N126").Value = -19334

# --- CRP row 134 ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 2837.88
$ws.Range("I134").Value = 1867
$ws.Range("K134").Value = 5601
$ws.Range("M134").Value = -3066

# --- CUL row 113 ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 2854.1785
$ws.Range("I113").Value = 5000
$ws.Range("K113").Value = 15000
$ws.Range("M113").Value = -12830

# --- CUL row 134 ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H134").Value = 3500
$ws.Range("I134").Value = 3500
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 10500
$ws.Range("L134").Value = 0
$ws.Range("M134").Value = -5430
$ws.Range("N134").ClearContents()

# --- CUL row 140 ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H140").Value = 2381.4443
$ws.Range("I140").Value = 2179.1875
$ws.Range("K140").Value = 6537.5625
$ws.Range("M140").Value = -1357.5625

# --- GSM row 70 ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 11416.5
$ws.Range("J70").Value = 5999
$ws.Range("L70").Value = 5999
$ws.Range("N70").Value = -6539

# --- GSM row 73 ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H73").Value = 11416.5
$ws.Range("J73").Value = 5999
$ws.Range("L73").Value = 5999
$ws.Range("N73").Value = -7871

# --- GSM row 94 ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H94").Value = 48905
$ws.Range("I94").Value = 48905
$ws.Range("K94").Value = 48905
$ws.Range("M94").Value = -48229

# --- GSM row 97 ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 994.45
$ws.Range("I97").Value = 1207.7333
$ws.Range("J97").Value = 354.6
$ws.Range("K97").Value = 1207.7333
$ws.Range("L97").Value = 354.6
$ws.Range("M97").Value = -711.7333000000001
$ws.Range("N97").Value = -1346.6

# --- GSM row 122 ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 291531.72
$ws.Range("I122").Value = 471397.22
$ws.Range("J122").Value = 6744.6665
$ws.Range("K122").Value = 1414191.66
$ws.Range("L122").Value = 20233.9995
$ws.Range("M122").Value = -1411741.66
$ws.Range("N122").Value = -25133.9995

# --- GSM row 126 ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 3557.5881
$ws.Range("I126").Value = 3320.3572
$ws.Range("K126").Value = 9961.071599999999
$ws.Range("M126").Value = -7491.071599999999

# --- LTW row 45 ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H45").Value = 21499.666
$ws.Range("I45").Value = 21499.666
$ws.Range("K45").Value = 21499.666
$ws.Range("M45").Value = -21092.666

# --- LTW row 93 ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 2673.8462
$ws.Range("I93").Value = 2503.389
$ws.Range("J93").Value = 3057.375
$ws.Range("K93").Value = 2503.389
$ws.Range("L93").Value = 3057.375
$ws.Range("M93").Value = -1255.389
$ws.Range("N93").Value = -5553.375

# --- WVR row 62 ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 9718.333000000001
$ws.Range("J62").Value = 9718.333000000001
$ws.Range("L62").Value = 9718.333000000001
$ws.Range("N62").Value = -10966.333

# --- WVR row 65 ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H65").Value = 9718.333000000001
$ws.Range("J65").Value = 9718.333000000001
$ws.Range("L65").Value = 48591.665
$ws.Range("N65").Value = -54831.665

# --- WVR row 81 ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 1458.3334
$ws.Range("I81").Value = 1450
$ws.Range("J81").Value = 1500
$ws.Range("K81").Value = 2900
$ws.Range("L81").Value = 3000
$ws.Range("M81").Value = -1839
$ws.Range("N81").Value = -5122

# --- WVR row 84 ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H84").Value = 1458.3334
$ws.Range("I84").Value = 1450
$ws.Range("J84").Value = 1500
$ws.Range("K84").Value = 14500
$ws.Range("L84").Value = 15000
$ws.Range("M84").Value = -9196
$ws.Range("N84").Value = -25608

# --- WVR row 113 ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 898.3333
$ws.Range("I113").Value = 650.0769
$ws.Range("J113").Value = 1301.75
$ws.Range("K113").Value = 1950.2307
$ws.Range("L113").Value = 3905.25
$ws.Range("M113").Value = 219.7692999999999
$ws.Range("N113").Value = -8245.25

# --- WVR row 122 ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 2506.111
$ws.Range("I122").Value = 1201.6154
$ws.Range("K122").Value = 3604.8462
$ws.Range("M122").Value = -1154.8462

# --- WVR row 136 ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 3102.1035
$ws.Range("I136").Value = 2426.7144
$ws.Range("K136").Value = 7280.1432
$ws.Range("M136").Value = -4730.1432
